$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells C1:H1 (bold, centered, top-aligned, thin border - matching A1/B1 style) ---
$headerCells = @("C1","D1","E1","F1","G1","H1")
foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box all sides)
}

$ws.Range("C1").Value = "Título"
$ws.Range("D1").Value = "Ubicación"
$ws.Range("E1").Value = "Paga"
$ws.Range("F1").Value = "Empresa"
$ws.Range("G1").Value = "I_Detalles"
$ws.Range("H1").Value = "I_Descripción"

# --- Row 2: B2 URL changes to the new job id; C2:H2 become empty (but present) string cells ---
$ws.Range("B2").Value = "https://www.occ.com.mx/empleos/de-especialista-en-implementacion-de-software/?jobid=20642450"
foreach ($addr in @("C2","D2","E2","F2","G2","H2")) {
    # Touch the cell so it becomes a real (empty) entry in the sheet and the
    # used range grows to include it, matching the <c .. t="inlineStr"/> cells
    # that appear (empty) in row 2 of the target.
    $ws.Range($addr).NumberFormat = "General"
}

# --- Row 3: A3/B3 updated, C3:H3 filled with the scraped job detail ---
$ws.Range("A3").Value = "Postularme`nEsta es una vacante externa, deberás completar el proceso en el sitio de la empresa.`nPostularme`nPostularme`nPostularme"
$ws.Range("B3").Value = "https://www.occ.com.mx/empleos/de-especialista-en-implementacion-de-software/?jobid=20652795"
$ws.Range("C3").Value = "Especialista Técnico en Integración/Implementación de Sistemas"
$ws.Range("D3").Value = "Ciudad de México"
$ws.Range("E3").Value = "Si el reclutador te contacta podrás conocer el sueldo"
$ws.Range("F3").Value = "Weldmation de México, S.A. de C.V.  en`nCiudad de México"
$ws.Range("G3").Value = "Detalles`nContratación:`nTiempo completo`nHorario:`nJornada completa`nEspacio de trabajo:`n`t`t`t`t`t`tPresencial"
$ws.Range("H3").Value = "Descripción`nRequerimientos`nEducación mínima: Diplomado`nObjetivo del Rol`nActuar como puente técnico entre el equipo interno de TI, el área de negocio y los equipos de desarrollo externos (DMS) facilitando la implementación eficiente de integraciones tecnológicas, garantizando seguridad, rendimiento y estabilidad en los sistemas.`n Responsabilidades Técnicas?`n· `nA`nnalizar, diseñar e implementar integraciones entre sistemas utilizando APIs (REST/JSON/XML) y protocolos de envío como SFTP.`n· Identificar y resolver proactivamente problemas técnicos relacionados con la configuración de integraciones e implementaciones.`n· Asegurar el cumplimiento de estándares de seguridad y rendimiento en las integraciones implementadas.`n· Participar activamente en las etapas del ciclo de vida del desarrollo de software, especialmente durante procesos de implementación y despliegue.`n· Colaborar con equipos técnicos de los distribuidores para definir soluciones de integración alineadas a los requerimientos del negocio.`n· Generar documentación técnica clara y precisa para facilitar la adopción y operación de las integraciones.`n· Brindar capacitación técnica a propietarios de aplicaciones del negocio y TI`n en el `nuso de herramientas y prácticas relacionadas con las integraciones.`n· Mantener comunicación directa con los equipos de TI de DTNA y DMS para resolución de incidencias y mejora continua.`n· Aplicar prácticas de trabajo bajo metodologías ágiles como Scrum para el seguimiento técnico de los desarrollos (deseable).`nConocimientos y Herramientas Deseadas`n· API Management, Web Services, JSON, REST, XML.`n· Protocolos de comunicación seguros como SFTP.`n· Seguridad de integraciones (tokens, autenticación, certificados).`n· Control de versiones y herramientas de deployment (Git, CI/CD).`n· Experiencia en ambientes con múltiples sistemas y plataformas (DEV/QA/PRD).`n· Comprensión de metodologías ágiles (Scrum) y sus herramientas (Jira, Confluence, etc.).`n· Conocimiento amplio de herramientas de observabilidad y monitoreo (Splunk, ThosandEyes, Etc.)`n· Habilidad para entender procesos de negocio y transformarlos en un requerimiento técnico`n· Habilidades de comunicación efectiva`n· Experiencia en Industria Automotriz deseable.`n· Inglés mínimo B2"

Write-Host "edit applied"
